$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "default payment prediction" (L) values that flipped between 0 and 1
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(27, 12).Value = 1
$ws.Cells.Item(30, 12).Value = 1
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(55, 12).Value = 1
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(70, 12).Value = 1
$ws.Cells.Item(71, 12).Value = 1
$ws.Cells.Item(75, 12).Value = 1
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(93, 12).Value = 1
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(111, 12).Value = 0

# Update the "probability" (M) text values (kept as text, not numbers, via a
# quoted-formula + paste-values round-trip so the cell stays a text/string cell
# instead of being coerced into a numeric cell by the normal Value setter).
$ws.Cells.Item(8, 13).Formula = "=""0.75"""
$ws.Cells.Item(10, 13).Formula = "=""0.67"""
$ws.Cells.Item(12, 13).Formula = "=""0.91"""
$ws.Cells.Item(13, 13).Formula = "=""0.98"""
$ws.Cells.Item(14, 13).Formula = "=""0.93"""
$ws.Cells.Item(15, 13).Formula = "=""0.92"""
$ws.Cells.Item(16, 13).Formula = "=""0.91"""
$ws.Cells.Item(17, 13).Formula = "=""0.73"""
$ws.Cells.Item(18, 13).Formula = "=""0.81"""
$ws.Cells.Item(19, 13).Formula = "=""0.76"""
$ws.Cells.Item(20, 13).Formula = "=""0.64"""
$ws.Cells.Item(21, 13).Formula = "=""0.61"""
$ws.Cells.Item(22, 13).Formula = "=""0.66"""
$ws.Cells.Item(23, 13).Formula = "=""0.83"""
$ws.Cells.Item(24, 13).Formula = "=""0.96"""
$ws.Cells.Item(25, 13).Formula = "=""0.75"""
$ws.Cells.Item(26, 13).Formula = "=""0.61"""
$ws.Cells.Item(27, 13).Formula = "=""0.54"""
$ws.Cells.Item(28, 13).Formula = "=""0.69"""
$ws.Cells.Item(29, 13).Formula = "=""0.84"""
$ws.Cells.Item(30, 13).Formula = "=""0.56"""
$ws.Cells.Item(31, 13).Formula = "=""0.76"""
$ws.Cells.Item(32, 13).Formula = "=""0.86"""
$ws.Cells.Item(33, 13).Formula = "=""0.86"""
$ws.Cells.Item(34, 13).Formula = "=""0.82"""
$ws.Cells.Item(35, 13).Formula = "=""0.93"""
$ws.Cells.Item(36, 13).Formula = "=""0.68"""
$ws.Cells.Item(37, 13).Formula = "=""0.63"""
$ws.Cells.Item(38, 13).Formula = "=""0.85"""
$ws.Cells.Item(39, 13).Formula = "=""0.74"""
$ws.Cells.Item(40, 13).Formula = "=""0.89"""
$ws.Cells.Item(41, 13).Formula = "=""0.7"""
$ws.Cells.Item(43, 13).Formula = "=""0.9"""
$ws.Cells.Item(44, 13).Formula = "=""0.66"""
$ws.Cells.Item(45, 13).Formula = "=""0.93"""
$ws.Cells.Item(46, 13).Formula = "=""0.61"""
$ws.Cells.Item(48, 13).Formula = "=""0.67"""
$ws.Cells.Item(50, 13).Formula = "=""0.97"""
$ws.Cells.Item(51, 13).Formula = "=""0.7"""
$ws.Cells.Item(52, 13).Formula = "=""0.85"""
$ws.Cells.Item(53, 13).Formula = "=""0.89"""
$ws.Cells.Item(54, 13).Formula = "=""0.85"""
$ws.Cells.Item(55, 13).Formula = "=""0.61"""
$ws.Cells.Item(56, 13).Formula = "=""0.69"""
$ws.Cells.Item(57, 13).Formula = "=""0.75"""
$ws.Cells.Item(58, 13).Formula = "=""0.93"""
$ws.Cells.Item(59, 13).Formula = "=""0.9"""
$ws.Cells.Item(60, 13).Formula = "=""0.58"""
$ws.Cells.Item(61, 13).Formula = "=""0.86"""
$ws.Cells.Item(62, 13).Formula = "=""0.72"""
$ws.Cells.Item(63, 13).Formula = "=""0.51"""
$ws.Cells.Item(64, 13).Formula = "=""0.99"""
$ws.Cells.Item(65, 13).Formula = "=""0.7"""
$ws.Cells.Item(66, 13).Formula = "=""0.52"""
$ws.Cells.Item(67, 13).Formula = "=""0.87"""
$ws.Cells.Item(68, 13).Formula = "=""0.9"""
$ws.Cells.Item(69, 13).Formula = "=""0.54"""
$ws.Cells.Item(70, 13).Formula = "=""0.53"""
$ws.Cells.Item(71, 13).Formula = "=""0.84"""
$ws.Cells.Item(72, 13).Formula = "=""0.57"""
$ws.Cells.Item(73, 13).Formula = "=""0.96"""
$ws.Cells.Item(74, 13).Formula = "=""0.71"""
$ws.Cells.Item(75, 13).Formula = "=""0.58"""
$ws.Cells.Item(76, 13).Formula = "=""0.9"""
$ws.Cells.Item(77, 13).Formula = "=""0.66"""
$ws.Cells.Item(78, 13).Formula = "=""0.58"""
$ws.Cells.Item(79, 13).Formula = "=""0.84"""
$ws.Cells.Item(80, 13).Formula = "=""0.74"""
$ws.Cells.Item(81, 13).Formula = "=""0.77"""
$ws.Cells.Item(82, 13).Formula = "=""0.76"""
$ws.Cells.Item(83, 13).Formula = "=""0.92"""
$ws.Cells.Item(84, 13).Formula = "=""0.57"""
$ws.Cells.Item(85, 13).Formula = "=""0.69"""
$ws.Cells.Item(86, 13).Formula = "=""0.82"""
$ws.Cells.Item(87, 13).Formula = "=""0.87"""
$ws.Cells.Item(88, 13).Formula = "=""0.88"""
$ws.Cells.Item(89, 13).Formula = "=""0.62"""
$ws.Cells.Item(91, 13).Formula = "=""0.99"""
$ws.Cells.Item(93, 13).Formula = "=""0.62"""
$ws.Cells.Item(94, 13).Formula = "=""0.9"""
$ws.Cells.Item(95, 13).Formula = "=""0.63"""
$ws.Cells.Item(96, 13).Formula = "=""0.65"""
$ws.Cells.Item(97, 13).Formula = "=""0.62"""
$ws.Cells.Item(98, 13).Formula = "=""0.77"""
$ws.Cells.Item(99, 13).Formula = "=""0.77"""
$ws.Cells.Item(101, 13).Formula = "=""0.63"""
$ws.Cells.Item(102, 13).Formula = "=""0.82"""
$ws.Cells.Item(103, 13).Formula = "=""0.94"""
$ws.Cells.Item(104, 13).Formula = "=""0.62"""
$ws.Cells.Item(105, 13).Formula = "=""0.78"""
$ws.Cells.Item(106, 13).Formula = "=""0.94"""
$ws.Cells.Item(107, 13).Formula = "=""0.76"""
$ws.Cells.Item(108, 13).Formula = "=""0.55"""
$ws.Cells.Item(109, 13).Formula = "=""0.66"""
$ws.Cells.Item(110, 13).Formula = "=""0.7"""
$ws.Cells.Item(111, 13).Formula = "=""0.72"""

$mRange = $ws.Range("M2:M111")
$mRange.Copy()
$mRange.PasteSpecial(-4163)
